# The workbook has a 5-column table (A:E, headers A,B,C,D,F) in A1:E25.
# This edit inserts a brand-new leading "ID" column so the table becomes
# A:F (headers ID,A,B,C,D,F) in A1:F25, shifting the original data from
# A:E to B:F and populating the new column A with row identifiers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing A:E columns one place to the right (-> B:F), pushing
# every existing cell (values, types and styles) over intact.
$ws.Columns.Item(1).Insert()

# Give the new header cell (A1) the same look as the rest of the header
# row (bold, centered, thin border) by cloning the format from an
# existing header cell, then set its text.
$ws.Range("F1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("A1").Value2 = "ID"

# Fill in the new identifier column for each data row (2-25).
$ws.Cells.Item(2, 1).Value2 = "Hb 2"
$ws.Cells.Item(3, 1).Value2 = "Hb 3"
$ws.Cells.Item(4, 1).Value2 = "S 24"
$ws.Cells.Item(5, 1).Value2 = "S 28"
$ws.Cells.Item(6, 1).Value2 = "Hb 107"
$ws.Cells.Item(7, 1).Value2 = "Hb 66"
$ws.Cells.Item(8, 1).Value2 = "Hb 69"
$ws.Cells.Item(9, 1).Value2 = "Hb 95"
$ws.Cells.Item(10, 1).Value2 = "Hb 99"
$ws.Cells.Item(11, 1).Value2 = "Hb 92"
$ws.Cells.Item(12, 1).Value2 = "Hb 40"
$ws.Cells.Item(13, 1).Value2 = "Hb 41"
$ws.Cells.Item(14, 1).Value2 = "S 11"
$ws.Cells.Item(15, 1).Value2 = "Hb 57"
$ws.Cells.Item(16, 1).Value2 = "S 21"
$ws.Cells.Item(17, 1).Value2 = "S 22"
$ws.Cells.Item(18, 1).Value2 = "S 3"
$ws.Cells.Item(19, 1).Value2 = "S 4"
$ws.Cells.Item(20, 1).Value2 = "S 5"
$ws.Cells.Item(21, 1).Value2 = "Hb 74"
$ws.Cells.Item(22, 1).Value2 = "Hb 79"
$ws.Cells.Item(23, 1).Value2 = "Hb 32"
$ws.Cells.Item(24, 1).Value2 = "S 15"
$ws.Cells.Item(25, 1).Value2 = "S 16"

[void]$ws.Range("A1").Select()
